$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FortiBranch")
$ws.Activate()

# Update the Destination column (B) for rows 3-5 to the new value "appdb_Server"
$ws.Range("B3").Value = "appdb_Server"
$ws.Range("B4").Value = "appdb_Server"
$ws.Range("B5").Value = "appdb_Server"

# Update the selected cell to match the diff (B5 selected instead of E10)
$ws.Range("B5").Select()
